$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1265.8077
$ws.Range("I137").Value = 908.3889
$ws.Range("J137").Value = 2070
$ws.Range("K137").Value = 2725.1667
$ws.Range("L137").Value = 6210
$ws.Range("M137").Value = -175.1667000000002
$ws.Range("N137").Value = -11310

$ws.Range("H138").Value = 501503.03
$ws.Range("I138").Value = 1649.2222
$ws.Range("J138").Value = 631899.7
$ws.Range("K138").Value = 4947.6666
$ws.Range("L138").Value = 1895699.1
$ws.Range("M138").Value = 192.3334000000004
$ws.Range("N138").Value = -1905979.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 14880.143
$ws.Range("I2").Value = 469.6
$ws.Range("J2").Value = 50906.5
$ws.Range("K2").Value = 469.6
$ws.Range("L2").Value = 50906.5
$ws.Range("M2").Value = -356.6
$ws.Range("N2").Value = -51132.5

$ws.Range("H5").Value = 190.6
$ws.Range("J5").Value = 200
$ws.Range("L5").Value = 200
$ws.Range("N5").Value = -424

$ws.Range("H32").Value = 3375.5833
$ws.Range("I32").Value = 3365.2
$ws.Range("J32").Value = 3427.5
$ws.Range("K32").Value = 3365.2
$ws.Range("L32").Value = 3427.5
$ws.Range("M32").Value = -3078.2
$ws.Range("N32").Value = -4001.5

$ws.Range("H61").Value = 1894.2222
$ws.Range("I61").Value = 1512.5
$ws.Range("J61").Value = 2199.6
$ws.Range("K61").Value = 1512.5
$ws.Range("L61").Value = 2199.6
$ws.Range("M61").Value = -1300.5
$ws.Range("N61").Value = -2623.6

$ws.Range("H74").Value = 1725
$ws.Range("I74").Value = 839.3333
$ws.Range("J74").Value = 2389.25
$ws.Range("K74").Value = 839.3333
$ws.Range("L74").Value = 2389.25
$ws.Range("M74").Value = 34.66669999999999
$ws.Range("N74").Value = -4137.25

$ws.Range("H77").Value = 1725
$ws.Range("I77").Value = 839.3333
$ws.Range("J77").Value = 2389.25
$ws.Range("K77").Value = 4196.6665
$ws.Range("L77").Value = 11946.25
$ws.Range("M77").Value = 171.3334999999997
$ws.Range("N77").Value = -20682.25

$ws.Range("H116").Value = 14880.143
$ws.Range("I116").Value = 469.6
$ws.Range("J116").Value = 50906.5
$ws.Range("K116").Value = 469.6
$ws.Range("L116").Value = 50906.5
$ws.Range("M116").Value = 1824.4
$ws.Range("N116").Value = -55494.5

$ws.Range("H132").Value = 1943.325
$ws.Range("I132").Value = 1663.1562
$ws.Range("J132").Value = 3064
$ws.Range("K132").Value = 4989.4686
$ws.Range("L132").Value = 9192
$ws.Range("M132").Value = -2459.4686
$ws.Range("N132").Value = -14252

$ws.Range("H136").Value = 1894.2222
$ws.Range("I136").Value = 1512.5
$ws.Range("J136").Value = 2199.6
$ws.Range("K136").Value = 4537.5
$ws.Range("L136").Value = 6598.799999999999
$ws.Range("M136").Value = -1987.5
$ws.Range("N136").Value = -11698.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 14880.143
$ws.Range("I3").Value = 469.6
$ws.Range("J3").Value = 50906.5
$ws.Range("K3").Value = 469.6
$ws.Range("L3").Value = 50906.5
$ws.Range("M3").Value = -355.6
$ws.Range("N3").Value = -51134.5

$ws.Range("H4").Value = 190.6
$ws.Range("J4").Value = 200
$ws.Range("L4").Value = 200
$ws.Range("N4").Value = -430

$ws.Range("H105").Value = 142860140
$ws.Range("I105").Value = 142860140
$ws.Range("K105").Value = 142860140
$ws.Range("M105").Value = -142858393

$ws.Range("H134").Value = 5324.148
$ws.Range("I134").Value = 1013.13043
$ws.Range("J134").Value = 30112.5
$ws.Range("K134").Value = 3039.39129
$ws.Range("L134").Value = 90337.5
$ws.Range("M134").Value = -504.39129
$ws.Range("N134").Value = -95407.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2082.0417
$ws.Range("I31").Value = 1902.5
$ws.Range("K31").Value = 1902.5
$ws.Range("M31").Value = -1607.5

$ws.Range("H34").Value = 2082.0417
$ws.Range("I34").Value = 1902.5
$ws.Range("K34").Value = 1902.5
$ws.Range("M34").Value = -1700.5

$ws.Range("H58").Value = 1447.3889
$ws.Range("I58").Value = 1088.1666
$ws.Range("J58").Value = 2165.8333
$ws.Range("K58").Value = 1088.1666
$ws.Range("L58").Value = 2165.8333
$ws.Range("M58").Value = -885.1666
$ws.Range("N58").Value = -2571.8333

$ws.Range("H132").Value = 10393.429
$ws.Range("I132").Value = 16970.857
$ws.Range("K132").Value = 50912.571
$ws.Range("M132").Value = -48382.571

$ws.Range("H134").Value = 2570.1667
$ws.Range("I134").Value = 2816.5833
$ws.Range("J134").Value = 2077.3333
$ws.Range("K134").Value = 8449.749899999999
$ws.Range("L134").Value = 6231.999899999999
$ws.Range("M134").Value = -5914.749899999999
$ws.Range("N134").Value = -11301.9999

$ws.Range("H136").Value = 1447.3889
$ws.Range("I136").Value = 1088.1666
$ws.Range("J136").Value = 2165.8333
$ws.Range("K136").Value = 3264.4998
$ws.Range("L136").Value = 6497.499899999999
$ws.Range("M136").Value = -714.4998000000001
$ws.Range("N136").Value = -11597.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 34619052
$ws.Range("I70").Value = 27781566
$ws.Range("J70").Value = 50003400
$ws.Range("K70").Value = 27781566
$ws.Range("L70").Value = 50003400
$ws.Range("M70").Value = -27781296
$ws.Range("N70").Value = -50003940

$ws.Range("H73").Value = 34619052
$ws.Range("I73").Value = 27781566
$ws.Range("J73").Value = 50003400
$ws.Range("K73").Value = 27781566
$ws.Range("L73").Value = 50003400
$ws.Range("M73").Value = -27780630
$ws.Range("N73").Value = -50005272

$ws.Range("H80").Value = 4119
$ws.Range("I80").Value = 3448.75
$ws.Range("K80").Value = 3448.75
$ws.Range("M80").Value = -2450.75

$ws.Range("H83").Value = 4119
$ws.Range("I83").Value = 3448.75
$ws.Range("K83").Value = 17243.75
$ws.Range("M83").Value = -12251.75

$ws.Range("H92").Value = 29999.75
$ws.Range("J92").Value = 29999.75
$ws.Range("L92").Value = 29999.75
$ws.Range("N92").Value = -33743.75

$ws.Range("H109").Value = 29999.666
$ws.Range("J109").Value = 29999.666
$ws.Range("L109").Value = 29999.666
$ws.Range("N109").Value = -32079.666

$ws.Range("H122").Value = 1368.3871
$ws.Range("I122").Value = 1483.174
$ws.Range("J122").Value = 1038.375
$ws.Range("K122").Value = 4449.522
$ws.Range("L122").Value = 3115.125
$ws.Range("M122").Value = -1999.522
$ws.Range("N122").Value = -8015.125

$ws.Range("H126").Value = 1991.7142
$ws.Range("I126").Value = 1658
$ws.Range("J126").Value = 2659.1428
$ws.Range("K126").Value = 4974
$ws.Range("L126").Value = 7977.428400000001
$ws.Range("M126").Value = -2504
$ws.Range("N126").Value = -12917.4284

$ws.Range("H132").Value = 2204.6155
$ws.Range("I132").Value = 2006.7778
$ws.Range("J132").Value = 2649.75
$ws.Range("K132").Value = 6020.3334
$ws.Range("L132").Value = 7949.25
$ws.Range("M132").Value = -3490.3334
$ws.Range("N132").Value = -13009.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1810
$ws.Range("J7").Value = 2050
$ws.Range("L7").Value = 2050
$ws.Range("N7").Value = -2274

$ws.Range("H126").Value = 1810
$ws.Range("J126").Value = 2050
$ws.Range("L126").Value = 6150
$ws.Range("N126").Value = -11090

$ws.Range("H132").Value = 22562.541
$ws.Range("I132").Value = 1314.8
$ws.Range("K132").Value = 3944.4
$ws.Range("M132").Value = -1414.4

$ws.Range("H136").Value = 5938.6
$ws.Range("I136").Value = 7173.25
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 21519.75
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -18969.75
$ws.Range("N136").Value = -8100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 23638556
$ws.Range("I122").Value = 23638556
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 70915668
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -70913218
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 3977.2856
$ws.Range("I132").Value = 10550
$ws.Range("J132").Value = 2881.8333
$ws.Range("K132").Value = 31650
$ws.Range("L132").Value = 8645.499899999999
$ws.Range("M132").Value = -29120
$ws.Range("N132").Value = -13705.4999

$ws.Range("H136").Value = 703.25
$ws.Range("I136").Value = 579.875
$ws.Range("J136").Value = 950
$ws.Range("K136").Value = 1739.625
$ws.Range("L136").Value = 2850
$ws.Range("M136").Value = 810.375
$ws.Range("N136").Value = -7950
